# Edit flujoefectivomensual1B.xlsx (Sheet1) to add an "ingreso" (income)
# section mirroring the existing "categoria" section, with its own totals
# row, and push the existing "subcat" (+ its TOTALES row) section further
# down the sheet.
#
# Strategy:
#  1. Copy/pin down cell FORMATS first (xlPasteFormats), referencing the
#     ORIGINAL (untouched) cells as sources, so that every destination cell
#     ends up re-using an existing style slot instead of Excel fabricating a
#     brand-new cellXfs entry per destination.
#  2. Only once every format-copy has been issued do we overwrite cell TEXT
#     VALUES - value writes never touch formatting, so this ordering makes
#     the whole thing independent of the fact that several rows are being
#     reshuffled (the source of a format-copy for a NEW row is often the
#     OLD row that will itself be overwritten with different text later).
#  3. Finally: clear the row that disappears completely (old row 16),
#     update the active-cell selection and the saved window position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Copy-Format($srcA1, $dstA1) {
    $ws.Range($srcA1).Copy()
    $ws.Range($dstA1).PasteSpecial(-4122)  # xlPasteFormats
}

# ---------------------------------------------------------------------
# 1) Format copies (sources are cells whose format is untouched at this
#    point in the script).
# ---------------------------------------------------------------------

# New row 11 ("TOTALES" under the new ingreso block): reuse style from the
# old row 14 TOTALES row (A14 plain style, B14/C14/D14 numeric style).
Copy-Format "A14" "A11"
Copy-Format "B14" "B11"
Copy-Format "B14" "C11"
Copy-Format "B14" "D11"

# New row 13 ("Concepto" label above the subcat block): reuse the style
# that used to live on row 16.
Copy-Format "A16" "A13"

# Row 14 becomes the subcat header row; A14 already has the right style
# (same plain bordered style both before and after), B14/C14/D14 need to
# take on the numeric styles that used to live on row 17.
Copy-Format "B17" "B14"
Copy-Format "C17" "C14"
Copy-Format "C17" "D14"

# Row 17 becomes the new subt TOTALES row; reuse the plain style that is
# already on A14 (same style index as A17 already has, and the one B/C/D17
# must switch to).
Copy-Format "A14" "B17"
Copy-Format "A14" "C17"
Copy-Format "A14" "D17"

# Row 8 (ingreso header row): B8/C8/D8 move from the numeric-protected
# style onto the same plain style A8 already uses.
Copy-Format "A8" "B8"
Copy-Format "A8" "C8"
Copy-Format "A8" "D8"

# ---------------------------------------------------------------------
# 2) Text values (safe to apply now; doesn't disturb formatting).
# ---------------------------------------------------------------------

# Row 8: "categoria" -> "ingreso"
$ws.Cells.Item(8, 1).Value2 = "{ingreso:nombre}"
$ws.Cells.Item(8, 2).Value2 = "{ingreso:banco1}"
$ws.Cells.Item(8, 3).Value2 = "{ingreso:total}"
$ws.Cells.Item(8, 4).Value2 = "{ingreso:pct}"

# Row 11 (new): totals for the ingreso block
$ws.Cells.Item(11, 1).Value2 = "TOTALES"
$ws.Cells.Item(11, 2).Value2 = "{ingresot:banco1}"
$ws.Cells.Item(11, 3).Value2 = "{ingresot:total}"
$ws.Cells.Item(11, 4).Value2 = "{ingresot:pct}"

# Row 13 (new): "Concepto" label, moved down from row 16
$ws.Cells.Item(13, 1).Value2 = "Concepto"

# Row 14: becomes the subcat header row (used to be the TOTALES row)
$ws.Cells.Item(14, 1).Value2 = "{subcat:nombre}"
$ws.Cells.Item(14, 2).Value2 = "{subcat:banco1}"
$ws.Cells.Item(14, 3).Value2 = "{subcat:total}"
$ws.Cells.Item(14, 4).Value2 = "{subcat:pct}"

# Row 17: becomes the subt TOTALES row (used to be the subcat header row)
$ws.Cells.Item(17, 1).Value2 = "TOTALES"
$ws.Cells.Item(17, 2).Value2 = "{subt:banco1}"
$ws.Cells.Item(17, 3).Value2 = "{subt:total}"
$ws.Cells.Item(17, 4).Value2 = "{subt:pct}"

# ---------------------------------------------------------------------
# 3) Row 16 disappears entirely now that its content lives on row 13.
# ---------------------------------------------------------------------
$ws.Cells.Item(16, 1).Clear()

# ---------------------------------------------------------------------
# 4) Selection + window position bookkeeping.
# ---------------------------------------------------------------------
$ws.Range("C18").Select()

try {
    $excel.ActiveWindow.Left = 2200
} catch {
}
